$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frecuencias")

# --- Convert relative-frequency fractions to rounded percentages (first table) ---
$ws.Range("C2").Value = 69.1
$ws.Range("C3").Value = 29.74
$ws.Range("C4").Value = 0.29
$ws.Range("C5").Value = 0.29
$ws.Range("C6").Value = 0.29
$ws.Range("C7").Value = 0.29

# --- Convert relative-frequency fractions to rounded percentages (second table) ---
$ws.Range("C11").Value = 38.48
$ws.Range("C12").Value = 27.7
$ws.Range("C13").Value = 14.87
$ws.Range("C14").Value = 11.66
$ws.Range("C15").Value = 7.29

# --- Insert a blank row + a "Total" row after the first table (rows 8 and 9) ---
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A9").Value = "Total"
$ws.Range("B9").Value = 343
$ws.Range("C9").Value = 100

# --- Insert a blank row + a "Total" row after the second table (now rows 18 and 19) ---
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(19).Insert()
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A19").Value = "Total"
$ws.Range("B19").Value = 343
$ws.Range("C19").Value = 100

# --- Remove the two table "title" cells that used to sit in column A of the header rows ---
$ws.Range("A1").Clear()
$ws.Range("A12").Clear()

# --- Reset column widths to default (drop the custom/bestFit widths) ---
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(2).ColumnWidth = 9.140625
$ws.Columns.Item(3).ColumnWidth = 9.140625
$ws.Columns.Item(1).UseStandardWidth = $true
$ws.Columns.Item(2).UseStandardWidth = $true
$ws.Columns.Item(3).UseStandardWidth = $true

# --- Clear the stale cell selection left over in the sheet view ---
$ws.Range("A1").Select()
